# Add a "Save" column (H) to the s_vals worksheet, matching the header
# style used by the other header cells (e.g. G1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "Save"

# Data values for H2:H14 as given by the diff
$values = @(0, 1, 0, 1, 0, 1, 0, 0, 0, 1, 0, 1, 1)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $values[$i]
}
